# "Generate Report for Archive"
# The localization status report is regenerated: the "Status" value that used
# to read "Ready for handoff" is now "In Translation" everywhere it appears
# (Overview!E2/F2 holding the per-language status, and the zh-cn/de-de detail
# sheets' Status cell, C2), and the now-narrower status column is shrunk to
# fit the shorter text on all three sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update every cell that shared the old "Ready for handoff" string so they
# all resolve to the same new shared string "In Translation".
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# The status columns are now narrower - shrink them to fit the new text.
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5
$zhcn.Columns("C").ColumnWidth = 12.5
$dede.Columns("C").ColumnWidth = 12.5
